# VRelayPaymentsCC.xlsx — refresh of recorded Katalon test-run timestamps/results
# (re-running the regression suite updates the "Date" column on each sheet's
# result row(s); CMCAutopayCC's run also flipped from Fail to Pass).

$wb = $excel.ActiveWorkbook

# PayNowCC — rows 2-7, column B (Date)
$ws = $wb.Worksheets.Item("PayNowCC")
$ws.Range("B2").Value = "Wed Mar 26 17:48:12 IST 2025"
$ws.Range("B3").Value = "Wed Mar 26 17:49:13 IST 2025"
$ws.Range("B4").Value = "Wed Mar 26 17:50:20 IST 2025"
$ws.Range("B5").Value = "Wed Mar 26 17:51:25 IST 2025"
$ws.Range("B6").Value = "Wed Mar 26 17:52:32 IST 2025"
$ws.Range("B7").Value = "Wed Mar 26 17:53:45 IST 2025"

# NoModifyAmountCC — row 2, column B (Date)
$ws = $wb.Worksheets.Item("NoModifyAmountCC")
$ws.Range("B2").Value = "Wed Mar 26 17:33:40 IST 2025"

# NoModifyBillingAddressCC — row 2, column B (Date)
$ws = $wb.Worksheets.Item("NoModifyBillingAddressCC")
$ws.Range("B2").Value = "Wed Mar 26 17:38:35 IST 2025"

# CCDeferredCC — row 2, column B (Date)
$ws = $wb.Worksheets.Item("CCDeferredCC")
$ws.Range("B2").Value = "Wed Mar 26 17:20:33 IST 2025"

# CMCAutopayCC — row 2: Result flips Fail -> Pass, and Date updates
$ws = $wb.Worksheets.Item("CMCAutopayCC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Mar 26 18:48:10 IST 2025"

# PayNowCreditCardDCF — row 2, column B (Date)
$ws = $wb.Worksheets.Item("PayNowCreditCardDCF")
$ws.Range("B2").Value = "Wed Mar 26 17:43:36 IST 2025"

# PayNowCreditCardSCF — row 2, column B (Date)
$ws = $wb.Worksheets.Item("PayNowCreditCardSCF")
$ws.Range("B2").Value = "Wed Mar 26 17:58:22 IST 2025"

# DCFCCVerbiage — row 2, column B (Date)
$ws = $wb.Worksheets.Item("DCFCCVerbiage")
$ws.Range("B2").Value = "Wed Mar 26 18:02:56 IST 2025"

# SCFCCVerbiage — row 2, column B (Date)
$ws = $wb.Worksheets.Item("SCFCCVerbiage")
$ws.Range("B2").Value = "Wed Mar 26 18:05:56 IST 2025"
